$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format before writing so numeric-looking price
# strings (e.g. "587.42") are stored as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '64.458.82'
$ws.Range('E2').Value = '  -0.24%  '

$ws.Range('D3').Value = '3.504.28'

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = '587.42'
$ws.Range('E5').Value = '  +0.32%  '

$ws.Range('D6').Value = '135.68'
$ws.Range('E6').Value = '  +2.55%  '

$ws.Range('D7').Value = '3.504.00'
$ws.Range('E7').Value = '  -0.07%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('E9').Value = '  -0.35%  '

$ws.Range('D10').Value = '0.124'
$ws.Range('E10').Value = '  +0.39%  '

$ws.Range('E11').Value = '  -0.98%  '

$ws.Range('E12').Value = '  -2.77%  '

$ws.Range('D13').Value = '4.102.38'
$ws.Range('E13').Value = '  -0.11%  '

$ws.Range('D14').Value = '0.0000181'
$ws.Range('E14').Value = '  +0.94%  '

$ws.Range('E15').Value = '  +1.31%  '

$ws.Range('D16').Value = '3.503.63'
$ws.Range('E16').Value = '  -0.18%  '

$ws.Range('D17').Value = '64.460.26'
$ws.Range('E17').Value = '  -0.22%  '

$ws.Range('D18').Value = '24.99'
$ws.Range('E18').Value = '  -9.63%  '

$ws.Range('D19').Value = '10.00'
$ws.Range('E19').Value = '  +0.78%  '

$ws.Range('D20').Value = '5.64'
$ws.Range('E20').Value = '  -0.83%  '

$ws.Range('D21').Value = '13.80'
$ws.Range('E21').Value = '  -2.75%  '

$ws.Range('D22').Value = '385.44'
$ws.Range('E22').Value = '  -1.48%  '

$ws.Range('E23').Value = '  -1.26%  '

$ws.Range('D24').Value = '3.645.23'
$ws.Range('E24').Value = '  -0.14%  '

$ws.Range('D25').Value = '74.12'
$ws.Range('E25').Value = '  +0.26%  '

$ws.Range('E26').Value = '  +0.05%  '

$ws.Range('E27').Value = '  +1.39%  '

$ws.Range('E28').Value = '  +3.74%  '

$ws.Range('E29').Value = '  -0.47%  '

$ws.Range('E30').Value = '  +0.75%  '

$ws.Range('E31').Value = '  +0.08%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '8.27'
$ws.Range('E32').Value = '  +0.70%  '

$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '2.23'
$ws.Range('E33').Value = '  -1.27%  '

$ws.Range('D34').Value = '3.524.95'
$ws.Range('E34').Value = '  +0.33%  '

$ws.Range('E35').Value = '  -0.02%  '

$ws.Range('E36').Value = '  +1.20%  '

$ws.Range('D37').Value = '23.61'
$ws.Range('E37').Value = '  -1.65%  '

$ws.Range('D38').Value = '5.33'
$ws.Range('E38').Value = '  +1.14%  '

$ws.Range('D39').Value = '1.55'
$ws.Range('E39').Value = '  -2.19%  '

$ws.Range('D40').Value = '6.85'
$ws.Range('E40').Value = '  -1.54%  '

$ws.Range('D41').Value = '163.50'
$ws.Range('E41').Value = '  -4.37%  '

$ws.Range('E42').Value = '  -2.58%  '

$ws.Range('D43').Value = '0.807'
$ws.Range('E43').Value = '  -0.75%  '

$ws.Range('D44').Value = '26.04'
$ws.Range('E44').Value = '  -1.43%  '

$ws.Range('E45').Value = '  +0.05%  '

$ws.Range('D46').Value = '41.88'
$ws.Range('E46').Value = '  -0.50%  '

$ws.Range('E47').Value = '  +0.24%  '

$ws.Range('E48').Value = '  +0.71%  '

$ws.Range('E49').Value = '  -0.01%  '

$ws.Range('D50').Value = '2.478.69'
$ws.Range('E50').Value = '  +1.03%  '

$ws.Range('E51').Value = '  -1.35%  '

# Restore default (unstyled) cell style now that the text values are set.
$ws.Range("D2:D51").Style = "Normal"
